# Applies the "cryptos list" refresh for Wed Jun 14 14:53:08 UTC 2023.
# Most rows only get updated Price (column D) and Volume(1h) (column E) figures.
# Rows 50/51 additionally swap places: "Decentraland" <-> "EnergySwap" change
# rank order, so their Coin name / Link / Price / Volume values are updated too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value. Values that look like plain numbers are
# prefixed with a leading apostrophe so Excel keeps them as text (matching the
# original inline-string cells, e.g. "1.0000" must not become the number 1).
$updates = @{
    "D2" = "26.000.32"
    "E2" = "  -0.02%  "
    "D3" = "1.741.98"
    "E3" = "  +0.08%  "
    "D4" = "'0.9998"
    "E4" = "  -0.12%  "
    "D5" = "'247.25"
    "E5" = "  +3.33%  "
    "D6" = "'1.000"
    "E6" = "  -0.13%  "
    "D7" = "'0.5045"
    "E7" = "  -4.52%  "
    "D8" = "'0.2747"
    "E8" = "  +0.71%  "
    "D9" = "'0.06188"
    "E9" = "  +0.48%  "
    "D10" = "1.749.58"
    "E10" = "  +0.54%  "
    "D11" = "'0.07267"
    "E11" = "  +1.21%  "
    "D12" = "'0.6550"
    "E12" = "  +2.26%  "
    "D13" = "'15.12"
    "E13" = "  +0.47%  "
    "D14" = "'4.680"
    "E14" = "  +1.62%  "
    "D15" = "'77.61"
    "E15" = "  +0.11%  "
    "D16" = "'1.0000"
    "E16" = "  -0.10%  "
    "D17" = "'0.9996"
    "E17" = "  -0.15%  "
    "D18" = "26.018.49"
    "E18" = "  +0.04%  "
    "D19" = "'11.93"
    "E19" = "  +1.41%  "
    "D20" = "'0.000006853"
    "E20" = "  +1.55%  "
    "D21" = "1.975.82"
    "E21" = "  +0.58%  "
    "D22" = "'4.474"
    "E22" = "  +2.86%  "
    "D23" = "'8.735"
    "E23" = "  +1.29%  "
    "D24" = "'5.399"
    "E24" = "  +2.80%  "
    "D25" = "'135.57"
    "E25" = "  -3.28%  "
    "D26" = "'1.506"
    "E26" = "  +0.26%  "
    "D27" = "'15.27"
    "E27" = "  +0.13%  "
    "D28" = "'1.787"
    "E28" = "  +1.35%  "
    "D29" = "'105.42"
    "E29" = "  -0.33%  "
    "D30" = "'3.946"
    "E30" = "  +3.00%  "
    "D31" = "'0.08173"
    "E31" = "  -2.59%  "
    "D32" = "'3.684"
    "E32" = "  +1.18%  "
    "D33" = "'0.04687"
    "E33" = "  +2.11%  "
    "D34" = "'2.656"
    "E34" = "  +0.00%  "
    "D35" = "'0.9975"
    "E35" = "  +0.51%  "
    "D36" = "'0.6132"
    "D37" = "'2.754"
    "E37" = "  +2.21%  "
    "D39" = "'1.928"
    "E39" = "  +0.09%  "
    "D40" = "'0.9998"
    "E40" = "  -0.19%  "
    "D41" = "'100.90"
    "E41" = "  +2.20%  "
    "E42" = "  +1.08%  "
    "D43" = "'0.7627"
    "E43" = "  +2.49%  "
    "E44" = "  +1.41%  "
    "D45" = "'0.1163"
    "E45" = "  +1.83%  "
    "D46" = "'6.315"
    "D47" = "'55.58"
    "E47" = "  +1.71%  "
    "E48" = "  -0.21%  "
    "D49" = "'30.73"
    "E49" = "  +0.07%  "
    "B50" = "EnergySwap"
    "C50" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D50" = "'7.635"
    "E50" = "  +1.66%  "
    "B51" = "Decentraland"
    "C51" = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
    "D51" = "'0.3474"
    "E51" = "  +0.93%  "
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
